# Natmi following Dr Hou advice
# Updates ligand/receptor expressing-cell counts (1 -> 3) and all
# downstream NATMI-derived expression/specificity metrics for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.70050833333333
$ws.Range("H2").Value = 62.101525
$ws.Range("I2").Value = 0.8277101186170105
$ws.Range("J2").Value = 0.8277101186170105
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 46.05975733333333
$ws.Range("N2").Value = 138.179272
$ws.Range("O2").Value = 0.8970651351272991
$ws.Range("P2").Value = 0.897065135127299
$ws.Range("Q2").Value = 953.4603905099776
$ws.Range("R2").Value = 8581.143514589799
$ws.Range("S2").Value = 0.7425098894034013
$ws.Range("T2").Value = 0.7425098894034012

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.70050833333333
$ws.Range("H3").Value = 62.101525
$ws.Range("I3").Value = 0.8277101186170105
$ws.Range("J3").Value = 0.8277101186170105
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.770761666666667
$ws.Range("N3").Value = 8.312284999999999
$ws.Range("O3").Value = 0.05396367312415441
$ws.Range("P3").Value = 0.0539636731241544
$ws.Range("Q3").Value = 57.35617497051388
$ws.Range("R3").Value = 516.2055747346249
$ws.Range("S3").Value = 0.04466627828260342
$ws.Range("T3").Value = 0.04466627828260342

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.70050833333333
$ws.Range("H4").Value = 62.101525
$ws.Range("I4").Value = 0.8277101186170105
$ws.Range("J4").Value = 0.8277101186170105
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.802173666666667
$ws.Range("N4").Value = 5.406521
$ws.Range("O4").Value = 0.03509934175535083
$ws.Range("P4").Value = 0.03509934175535083
$ws.Range("Q4").Value = 37.30591100494722
$ws.Range("R4").Value = 335.7531990445249
$ws.Range("S4").Value = 0.02905208032770042
$ws.Range("T4").Value = 0.02905208032770042

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.70050833333333
$ws.Range("H5").Value = 62.101525
$ws.Range("I5").Value = 0.8277101186170105
$ws.Range("J5").Value = 0.8277101186170105
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7122493333333333
$ws.Range("N5").Value = 2.136748
$ws.Range("O5").Value = 0.01387184999319569
$ws.Range("P5").Value = 0.01387184999319569
$ws.Range("Q5").Value = 14.74392326007777
$ws.Range("R5").Value = 132.6953093407
$ws.Range("S5").Value = 0.01148187060330538
$ws.Range("T5").Value = 0.01148187060330538

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7925996666666667
$ws.Range("H6").Value = 2.377799
$ws.Range("I6").Value = 0.03169210888681734
$ws.Range("J6").Value = 0.03169210888681734
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 46.05975733333333
$ws.Range("N6").Value = 138.179272
$ws.Range("O6").Value = 0.8970651351272991
$ws.Range("P6").Value = 0.897065135127299
$ws.Range("Q6").Value = 36.50694830914755
$ws.Range("R6").Value = 328.562534782328
$ws.Range("S6").Value = 0.02842988594102187
$ws.Range("T6").Value = 0.02842988594102187

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7925996666666667
$ws.Range("H7").Value = 2.377799
$ws.Range("I7").Value = 0.03169210888681734
$ws.Range("J7").Value = 0.03169210888681734
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.770761666666667
$ws.Range("N7").Value = 8.312284999999999
$ws.Range("O7").Value = 0.05396367312415441
$ws.Range("P7").Value = 0.0539636731241544
$ws.Range("Q7").Value = 2.196104773412778
$ws.Range("R7").Value = 19.764942960715
$ws.Range("S7").Value = 0.00171022260458332
$ws.Range("T7").Value = 0.00171022260458332

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.7925996666666667
$ws.Range("H8").Value = 2.377799
$ws.Range("I8").Value = 0.03169210888681734
$ws.Range("J8").Value = 0.03169210888681734
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.802173666666667
$ws.Range("N8").Value = 5.406521
$ws.Range("O8").Value = 0.03509934175535083
$ws.Range("P8").Value = 0.03509934175535083
$ws.Range("Q8").Value = 1.428402247475445
$ws.Range("R8").Value = 12.855620227279
$ws.Range("S8").Value = 0.001112372160766193
$ws.Range("T8").Value = 0.001112372160766193

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.7925996666666667
$ws.Range("H9").Value = 2.377799
$ws.Range("I9").Value = 0.03169210888681734
$ws.Range("J9").Value = 0.03169210888681734
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7122493333333333
$ws.Range("N9").Value = 2.136748
$ws.Range("O9").Value = 0.01387184999319569
$ws.Range("P9").Value = 0.01387184999319569
$ws.Range("Q9").Value = 0.5645285841835556
$ws.Range("R9").Value = 5.080757257651999
$ws.Range("S9").Value = 0.0004396281804459543
$ws.Range("T9").Value = 0.0004396281804459543

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.562510666666667
$ws.Range("H10").Value = 4.687532
$ws.Range("I10").Value = 0.06247701111592723
$ws.Range("J10").Value = 0.06247701111592723
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.05975733333333
$ws.Range("N10").Value = 138.179272
$ws.Range("O10").Value = 0.8970651351272991
$ws.Range("P10").Value = 0.897065135127299
$ws.Range("Q10").Value = 71.96886213741155
$ws.Range("R10").Value = 647.719759236704
$ws.Range("S10").Value = 0.05604594841905903
$ws.Range("T10").Value = 0.05604594841905902

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.562510666666667
$ws.Range("H11").Value = 4.687532
$ws.Range("I11").Value = 0.06247701111592723
$ws.Range("J11").Value = 0.06247701111592723
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.770761666666667
$ws.Range("N11").Value = 8.312284999999999
$ws.Range("O11").Value = 0.05396367312415441
$ws.Range("P11").Value = 0.0539636731241544
$ws.Range("Q11").Value = 4.329344658957777
$ws.Range("R11").Value = 38.96410193062
$ws.Range("S11").Value = 0.003371489005634058
$ws.Range("T11").Value = 0.003371489005634058

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.562510666666667
$ws.Range("H12").Value = 4.687532
$ws.Range("I12").Value = 0.06247701111592723
$ws.Range("J12").Value = 0.06247701111592723
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.802173666666667
$ws.Range("N12").Value = 5.406521
$ws.Range("O12").Value = 0.03509934175535083
$ws.Range("P12").Value = 0.03509934175535083
$ws.Range("Q12").Value = 2.815915577352444
$ws.Range("R12").Value = 25.343240196172
$ws.Range("S12").Value = 0.002192901965010782
$ws.Range("T12").Value = 0.002192901965010783

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.562510666666667
$ws.Range("H13").Value = 4.687532
$ws.Range("I13").Value = 0.06247701111592723
$ws.Range("J13").Value = 0.06247701111592723
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7122493333333333
$ws.Range("N13").Value = 2.136748
$ws.Range("O13").Value = 0.01387184999319569
$ws.Range("P13").Value = 0.01387184999319569
$ws.Range("Q13").Value = 1.112897180659555
$ws.Range("R13").Value = 10.016074625936
$ws.Range("S13").Value = 0.0008666717262233624
$ws.Range("T13").Value = 0.0008666717262233624

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.953751
$ws.Range("H14").Value = 5.861253
$ws.Range("I14").Value = 0.07812076138024482
$ws.Range("J14").Value = 0.07812076138024483
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 46.05975733333333
$ws.Range("N14").Value = 138.179272
$ws.Range("O14").Value = 0.8970651351272991
$ws.Range("P14").Value = 0.897065135127299
$ws.Range("Q14").Value = 89.98929694975732
$ws.Range("R14").Value = 809.903672547816
$ws.Range("S14").Value = 0.0700794113638168
$ws.Range("T14").Value = 0.07007941136381682

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.953751
$ws.Range("H15").Value = 5.861253
$ws.Range("I15").Value = 0.07812076138024482
$ws.Range("J15").Value = 0.07812076138024483
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.770761666666667
$ws.Range("N15").Value = 8.312284999999999
$ws.Range("O15").Value = 0.05396367312415441
$ws.Range("P15").Value = 0.0539636731241544
$ws.Range("Q15").Value = 5.413378377011666
$ws.Range("R15").Value = 48.72040539310499
$ws.Range("S15").Value = 0.004215683231333597
$ws.Range("T15").Value = 0.004215683231333597

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.953751
$ws.Range("H16").Value = 5.861253
$ws.Range("I16").Value = 0.07812076138024482
$ws.Range("J16").Value = 0.07812076138024483
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.802173666666667
$ws.Range("N16").Value = 5.406521
$ws.Range("O16").Value = 0.03509934175535083
$ws.Range("P16").Value = 0.03509934175535083
$ws.Range("Q16").Value = 3.520998603423666
$ws.Range("R16").Value = 31.688987430813
$ws.Range("S16").Value = 0.002741987301873426
$ws.Range("T16").Value = 0.002741987301873426

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.953751
$ws.Range("H17").Value = 5.861253
$ws.Range("I17").Value = 0.07812076138024482
$ws.Range("J17").Value = 0.07812076138024483
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.7122493333333333
$ws.Range("N17").Value = 2.136748
$ws.Range("O17").Value = 0.01387184999319569
$ws.Range("P17").Value = 0.01387184999319569
$ws.Range("Q17").Value = 1.391557847249333
$ws.Range("R17").Value = 12.524020625244
$ws.Range("S17").Value = 0.001083679483220992
$ws.Range("T17").Value = 0.001083679483220992
